$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 21:51"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Range("B4").Value = 3592043
$ws.Range("C4").Value = 46966
$ws.Range("D4").Value = 1624586
$ws.Range("E4").Value = 1827730
$ws.Range("G4").Value = 584
$ws.Range("H4").Value = 139727

# Row 6: India -> India
$ws.Range("B6").Value = 970169
$ws.Range("C6").Value = 32682
$ws.Range("D6").Value = 613735
$ws.Range("E6").Value = 331505
$ws.Range("G6").Value = 614
$ws.Range("H6").Value = 24929

# Row 11: España -> Sudafrica
$ws.Range("A11").Value = "Sudafrica"
$ws.Range("B11").Value = 311049
$ws.Range("C11").Value = 12757
$ws.Range("D11").Value = 160693
$ws.Range("E11").Value = 145903
$ws.Range("G11").Value = 107
$ws.Range("H11").Value = 4453

# Row 12: Sudafrica -> España
$ws.Range("A12").Value = "España"
$ws.Range("B12").Value = 304574
$ws.Range("C12").Value = 875
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 4
$ws.Range("H12").Value = 28413

# Row 19: Alemania -> Alemania
$ws.Range("B19").Value = 201050
$ws.Range("C19").Value = 284
$ws.Range("E19").Value = 5904

# Row 21: Francia -> Francia
$ws.Range("B21").Value = 173304
$ws.Range("C21").Value = 416
$ws.Range("D21").Value = 78820
$ws.Range("E21").Value = 64364
$ws.Range("G21").Value = 20
$ws.Range("H21").Value = 30120

# Row 26: Egipto -> Egipto
$ws.Range("B26").Value = 84843
$ws.Range("C26").Value = 913
$ws.Range("D26").Value = 26135
$ws.Range("E26").Value = 54641
$ws.Range("G26").Value = 59
$ws.Range("H26").Value = 4067

# Row 38: Ucrania -> Emiratos Arabes Unidos
$ws.Range("A38").Value = "Emiratos Arabes Unidos"
$ws.Range("B38").Value = 55848
$ws.Range("C38").Value = 275
$ws.Range("D38").Value = 46418
$ws.Range("E38").Value = 9095
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 335

# Row 39: Emiratos Arabes Unidos -> Ucrania
$ws.Range("A39").Value = "Ucrania"
$ws.Range("B39").Value = 55607
$ws.Range("C39").Value = 836
$ws.Range("D39").Value = 28131
$ws.Range("E39").Value = 26049
$ws.Range("G39").Value = 15
$ws.Range("H39").Value = 1427

# Row 67: Uzbekistan -> Uzbekistan
$ws.Range("E67").Value = 6052
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 71

# Row 79: Malasia -> Costa Rica
$ws.Range("A79").Value = "Costa Rica"
$ws.Range("B79").Value = 8986
$ws.Range("C79").Value = 504
$ws.Range("D79").Value = 2551
$ws.Range("E79").Value = 6395
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 40

# Row 80: Republica de Macedonia -> Malasia
$ws.Range("A80").Value = "Malasia"
$ws.Range("B80").Value = 8734
$ws.Range("C80").Value = 5
$ws.Range("D80").Value = 8526
$ws.Range("E80").Value = 86
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 122

# Row 81: Costa Rica -> Republica de Macedonia
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 8530
$ws.Range("C81").Value = 198
$ws.Range("D81").Value = 4565
$ws.Range("E81").Value = 3572
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 393

# Row 113: Mali -> Cuba
$ws.Range("A113").Value = "Cuba"
$ws.Range("B113").Value = 2438
$ws.Range("C113").Value = 6
$ws.Range("D113").Value = 2277
$ws.Range("E113").Value = 74
$ws.Range("H113").Value = 87

# Row 114: Cuba -> Mali
$ws.Range("A114").Value = "Mali"
$ws.Range("B114").Value = 2433
$ws.Range("C114").Value = 10
$ws.Range("D114").Value = 1764
$ws.Range("E114").Value = 548
$ws.Range("H114").Value = 121

# Row 126: Hong Kong -> Libia
$ws.Range("A126").Value = "Libia"
$ws.Range("C126").Value = 26
$ws.Range("D126").Value = 373
$ws.Range("E126").Value = 1173
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 43

# Row 127: Libia -> Hong Kong
$ws.Range("A127").Value = "Hong Kong"
$ws.Range("B127").Value = 1589
$ws.Range("C127").Value = 19
$ws.Range("D127").Value = 1241
$ws.Range("E127").Value = 338
$ws.Range("G127").Value = 2
$ws.Range("H127").Value = 10

# Row 129: Yemen -> Yemen
$ws.Range("B129").Value = 1526
$ws.Range("C129").Value = 10
$ws.Range("D129").Value = 694
$ws.Range("E129").Value = 399
$ws.Range("G129").Value = 4
$ws.Range("H129").Value = 433

# Row 157: Angola -> Angola
$ws.Range("B157").Value = 576
$ws.Range("C157").Value = 35
$ws.Range("D157").Value = 124
$ws.Range("E157").Value = 425
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 27

# Row 209: Islas Malvinas -> Groenlandia
$ws.Range("A209").Value = "Groenlandia"

# Row 210: Groenlandia -> Islas Malvinas
$ws.Range("A210").Value = "Islas Malvinas"
